$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API's V3")
$ws.Columns.Item(2).ColumnWidth = 8.140625
